# AFDP-3349: Add new access control rule (DocumentRepository - Default access)
# to the Assignment Rules table on Sheet1, row 28.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlEdgeTop = 8
$xlContinuous = 1
$xlLineStyleNone = -4142

# B28 - rule name (wrapped text, thin border all around)
$b28 = $ws.Range("B28")
$b28.Value = "DocumentRepository – Default access"
$b28.Borders.LineStyle = $xlContinuous
$b28.WrapText = $true

# C28 - object type (thin border, but no top edge - visually joins the row above)
$c28 = $ws.Range("C28")
$c28.Value = "DOC_REPO"
$c28.Borders.LineStyle = $xlContinuous
$c28.Borders.Item($xlEdgeTop).LineStyle = $xlLineStyleNone

# D28 - expression 1 (thin border all around, custom TRUE/FALSE display format)
$d28 = $ws.Range("D28")
$d28.Value = "participants.?[participantType == '*'].isEmpty()"
$d28.Borders.LineStyle = $xlContinuous
$d28.NumberFormat = """TRUE"";""TRUE"";""FALSE"""

# E28, F28 - expression 2/3, left blank, thin border all around
$ws.Range("E28").Borders.LineStyle = $xlContinuous
$ws.Range("F28").Borders.LineStyle = $xlContinuous

# G28 - literal participant assignment value, thin border all around
$g28 = $ws.Range("G28")
$g28.Value = "*, *"
$g28.Borders.LineStyle = $xlContinuous

# H28 - blank, thin border all around
$ws.Range("H28").Borders.LineStyle = $xlContinuous

$ws.Rows("28").RowHeight = 30

$ws.Range("B29").Select()
